$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 534. This shifts the existing rows
# 534:561 down to 535:562, preserving all their data/formatting
# (matches the diff: old row534 data now lives in row535, etc., and a
# brand-new record appears as the new row534; dimension grows to R562).
$ws.Rows("534:534").Insert()

# Populate the new row 534 with the new weekly record.
$ws.Range("A534").Value2 = 3
$ws.Range("B534").Value2 = "Femacal de La Calera"
$ws.Range("C534").Value2 = "Coquimbo"
$ws.Range("D534").Value2 = 45267
$ws.Range("E534").Value2 = 5
$ws.Range("F534").Value2 = 100112001
$ws.Range("G534").Value2 = "Berenjena"
$ws.Range("H534").Value2 = "Sin especificar"
$ws.Range("I534").Value2 = "Primera"
$ws.Range("J534").Value2 = 65
$ws.Range("K534").Value2 = 9000
$ws.Range("L534").Value2 = 9000
$ws.Range("M534").Value2 = 9000
$ws.Range("N534").Value2 = '$/caja 60 unidades'
$ws.Range("O534").Value2 = "Región de Arica y Parinacota"
$ws.Range("P534").Value2 = 150
$ws.Range("Q534").Value2 = 60
$ws.Range("R534").Value2 = "Hortaliza"
